$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1932.5714
$ws.Range("I137").Value = 1636
$ws.Range("J137").Value = 4750
$ws.Range("K137").Value = 4908
$ws.Range("L137").Value = 14250
$ws.Range("M137").Value = -2358
$ws.Range("N137").Value = -19350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5808.086
$ws.Range("I32").Value = 4342.636
$ws.Range("K32").Value = 4342.636
$ws.Range("M32").Value = -4055.636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 25558
$ws.Range("I45").Value = 25558
$ws.Range("K45").Value = 25558
$ws.Range("M45").Value = -25181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4525.0586
$ws.Range("I61").Value = 4809.242
$ws.Range("J61").Value = 4004.0557
$ws.Range("K61").Value = 4809.242
$ws.Range("L61").Value = 4004.0557
$ws.Range("M61").Value = -4597.242
$ws.Range("N61").Value = -4428.0557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2215.25
$ws.Range("I74").Value = 2105
$ws.Range("K74").Value = 2105
$ws.Range("M74").Value = -1231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2215.25
$ws.Range("I77").Value = 2105
$ws.Range("K77").Value = 10525
$ws.Range("M77").Value = -6157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1365.3704
$ws.Range("I110").Value = 1198.7368
$ws.Range("J110").Value = 1761.125
$ws.Range("K110").Value = 1198.7368
$ws.Range("L110").Value = 1761.125
$ws.Range("M110").Value = 846.2632000000001
$ws.Range("N110").Value = -5851.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 715603.8
$ws.Range("I122").Value = 830190.2
$ws.Range("J122").Value = 5168.4
$ws.Range("K122").Value = 2490570.6
$ws.Range("L122").Value = 15505.2
$ws.Range("M122").Value = -2488120.6
$ws.Range("N122").Value = -20405.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2780691.8
$ws.Range("I132").Value = 1715.0952
$ws.Range("J132").Value = 6671259
$ws.Range("K132").Value = 5145.2856
$ws.Range("L132").Value = 20013777
$ws.Range("M132").Value = -2615.2856
$ws.Range("N132").Value = -20018837

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 53119.168
$ws.Range("J135").Value = 53119.168
$ws.Range("L135").Value = 53119.168
$ws.Range("N135").Value = -63259.168

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4525.0586
$ws.Range("I136").Value = 4809.242
$ws.Range("J136").Value = 4004.0557
$ws.Range("K136").Value = 14427.726
$ws.Range("L136").Value = 12012.1671
$ws.Range("M136").Value = -11877.726
$ws.Range("N136").Value = -17112.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 32476.334
$ws.Range("J141").Value = 32476.334
$ws.Range("L141").Value = 32476.334
$ws.Range("N141").Value = -42836.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4672.641
$ws.Range("I134").Value = 5555.231
$ws.Range("K134").Value = 16665.693
$ws.Range("M134").Value = -14130.693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").ClearContents()
$ws.Range("H31").Value = 13470
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 13470
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 13470
$ws.Range("N31").Value = -14060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M34").ClearContents()
$ws.Range("H34").Value = 13470
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 13470
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 13470
$ws.Range("N34").Value = -13874

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10649.2
$ws.Range("I122").Value = 561
$ws.Range("K122").Value = 5049
$ws.Range("M122").Value = -2599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1622908.8
$ws.Range("I122").Value = 2948270.5
$ws.Range("J122").Value = 3022.2222
$ws.Range("K122").Value = 8844811.5
$ws.Range("L122").Value = 9066.6666
$ws.Range("M122").Value = -8842361.5
$ws.Range("N122").Value = -13966.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2831.587
$ws.Range("I132").Value = 3227.2727
$ws.Range("J132").Value = 2468.875
$ws.Range("K132").Value = 9681.8181
$ws.Range("L132").Value = 7406.625
$ws.Range("M132").Value = -7151.8181
$ws.Range("N132").Value = -12466.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1013350
$ws.Range("I82").Value = 1667658
$ws.Range("J82").Value = 228180.4
$ws.Range("K82").Value = 1667658
$ws.Range("L82").Value = 228180.4
$ws.Range("M82").Value = -1667297
$ws.Range("N82").Value = -228902.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1013350
$ws.Range("I85").Value = 1667658
$ws.Range("J85").Value = 228180.4
$ws.Range("K85").Value = 1667658
$ws.Range("L85").Value = 228180.4
$ws.Range("M85").Value = -1666410
$ws.Range("N85").Value = -230676.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 23820054
$ws.Range("I132").Value = 41682532
$ws.Range("J132").Value = 3416.5
$ws.Range("K132").Value = 125047596
$ws.Range("L132").Value = 10249.5
$ws.Range("M132").Value = -125045066
$ws.Range("N132").Value = -15309.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 41743
$ws.Range("J135").Value = 41743
$ws.Range("L135").Value = 41743
$ws.Range("N135").Value = -51883

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6851.628
$ws.Range("I136").Value = 4621.943
$ws.Range("J136").Value = 16606.5
$ws.Range("K136").Value = 13865.829
$ws.Range("L136").Value = 49819.5
$ws.Range("M136").Value = -11315.829
$ws.Range("N136").Value = -54919.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 67143.336
$ws.Range("J138").Value = 67143.336
$ws.Range("L138").Value = 67143.336
$ws.Range("N138").Value = -77423.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 45824
$ws.Range("J140").Value = 45824
$ws.Range("L140").Value = 45824
$ws.Range("N140").Value = -56184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1101
$ws.Range("I126").Value = 800.5
$ws.Range("J126").Value = 1501.6666
$ws.Range("K126").Value = 2401.5
$ws.Range("L126").Value = 4504.9998
$ws.Range("M126").Value = 68.5
$ws.Range("N126").Value = -9444.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1737.9773
$ws.Range("I132").Value = 1412.6666
$ws.Range("J132").Value = 2713.9092
$ws.Range("K132").Value = 4237.9998
$ws.Range("L132").Value = 8141.7276
$ws.Range("M132").Value = -1707.9998
$ws.Range("N132").Value = -13201.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 37928
$ws.Range("J133").Value = 37928
$ws.Range("L133").Value = 37928
$ws.Range("N133").Value = -48048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2840.5527
$ws.Range("I136").Value = 3674.2222
$ws.Range("K136").Value = 11022.6666
$ws.Range("M136").Value = -8472.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 55528.75
$ws.Range("J137").Value = 55528.75
$ws.Range("L137").Value = 55528.75
$ws.Range("N137").Value = -65728.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 69306
$ws.Range("J139").Value = 69306
$ws.Range("L139").Value = 69306
$ws.Range("N139").Value = -79586
